$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActorTable")

# Row 10 - Actor009 -> SciFiWarrior
$ws.Range("B10").Value = "CharName_SciFiWarrior"
$ws.Range("C10").Value = "CharDesc_SciFiWarrior"
$ws.Range("D10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("L10").Value = "SciFiWarrior"
$ws.Range("M10").Value = "Portrait_SciFiWarrior"
$ws.Range("N10").Value = 0.05

# Row 11 - Actor010 -> ChaosElemental
$ws.Range("B11").Value = "CharName_ChaosElemental"
$ws.Range("C11").Value = "CharDesc_ChaosElemental"
$ws.Range("I11").Value = 2
$ws.Range("L11").Value = "ChaosElemental"
$ws.Range("M11").Value = "Portrait_ChaosElemental"
$ws.Range("N11").Value = 0.05
$ws.Range("O11").Value = $true

# Row 12 - Actor011 -> SuperHero
$ws.Range("B12").Value = "CharName_SuperHero"
$ws.Range("C12").Value = "CharDesc_SuperHero"
$ws.Range("L12").Value = "SuperHero"
$ws.Range("M12").Value = "Portrait_SuperHero"
$ws.Range("N12").Value = 0.05

# Row 13 - Actor012 -> Meryl
$ws.Range("B13").Value = "CharName_Meryl"
$ws.Range("C13").Value = "CharDesc_Meryl"
$ws.Range("D13").Value = 1
$ws.Range("I13").Value = 0
$ws.Range("L13").Value = "Meryl"
$ws.Range("M13").Value = "Portrait_Meryl"
$ws.Range("N13").Value = 0.05

# Row 14 - Actor013 -> GreekWarrior
$ws.Range("B14").Value = "CharName_GreekWarrior"
$ws.Range("C14").Value = "CharDesc_GreekWarrior"
$ws.Range("D14").Value = 1
$ws.Range("L14").Value = "GreekWarrior"
$ws.Range("M14").Value = "Portrait_GreekWarrior"
$ws.Range("N14").Value = 0.05
